# Prepare submission: add a "word contribution" table to the summary sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("summary")
$ws.Select()

# Widen column B to fit the new word labels.
$ws.Columns.Item(2).ColumnWidth = 20.83203125

# Stray formatted cell left over between the two tables.
$ws.Range("H14").Borders.LineStyle = 1

$words = @(
    "perfect", "tears", "masterpiece", "innocence", "themes", "makes",
    "beauty", "highly", "enjoy", "small", "truly", "natural", "film",
    "story", "made"
)
$contributions = @(
    0.65901024454201296, 0.46467479083404500, 0.38625577955699902,
    0.37838579668583999, 0.34497046296608302, 0.29170352791486598,
    0.23911131171627600, 0.20765276989099901, 0.19800602792321800,
    0.16843686110262401, 0.13207275853949099, 0.07067555532488080,
    0.00430305645958309, -0.00650456910664081, -0.07051945637686580
)

# Header row.
$ws.Range("B16").Value = "word"
$ws.Range("C16").Value = "contribution"
$ws.Range("B16:C16").Borders.LineStyle = 1
$ws.Range("B16:C16").HorizontalAlignment = -4108

for ($i = 0; $i -lt $words.Length; $i++) {
    $row = 17 + $i
    $ws.Cells.Item($row, 2).Value = $words[$i]
    $ws.Cells.Item($row, 3).Value = $contributions[$i]
}

$dataRange = $ws.Range("B17:C31")
$dataRange.Borders.LineStyle = 1
$dataRange.HorizontalAlignment = -4108
$dataRange.Style = "Comma"

$ws.Range("B16").Select()
$ws.Range("B16:C31").Select()

$ws.PageSetup.Orientation = 1
